# Atualização de bases das ligas, do dia: 27-03-2024 às 20:23
#
# 1) Rows 104 and 107 (same kick-off date, 2023-11-11 09:00) swap places.
# 2) Rows 115 and 116 (same kick-off date, 2024-03-09 09:30) swap places.
# 3) Two new fixtures are appended as rows 122 and 123.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the data (columns B..AC) of rows 104 and 107 -------------------
$row104 = $ws.Range("B104:AC104").Value2
$row107 = $ws.Range("B107:AC107").Value2
$ws.Range("B104:AC104").Value = $row107
$ws.Range("B107:AC107").Value = $row104

# --- 2) Swap the data (columns B..AC) of rows 115 and 116 -------------------
$row115 = $ws.Range("B115:AC115").Value2
$row116 = $ws.Range("B116:AC116").Value2
$ws.Range("B115:AC115").Value = $row116
$ws.Range("B116:AC116").Value = $row115

# --- 3) Append the two new match rows (122 and 123) --------------------------
# Copy formatting (bold/border/centering on col A, date format on col E) from
# the last existing data row so the new rows look consistent.
$ws.Range("A121:E121").Copy()
$ws.Range("A122:E123").PasteSpecial(-4122)

# Row 122: Parnu JK Vaprus 3-2 FC Kuressaare
$ws.Cells.Item(122, 1).Value = 120
$ws.Cells.Item(122, 2).Value = 7721008
$ws.Cells.Item(122, 3).Value = "Estonia Meistriliiga"
$ws.Cells.Item(122, 4).Value = "Estonia Meistriliiga"
$ws.Cells.Item(122, 5).Value = 45368.3125
$ws.Cells.Item(122, 6).Value = "Parnu JK Vaprus"
$ws.Cells.Item(122, 7).Value = "FC Kuressaare"
$ws.Cells.Item(122, 8).Value = 3
$ws.Cells.Item(122, 9).Value = 2
$ws.Cells.Item(122, 10).Value = "H"
$row122vals = @(1.5, 4.5, 4.5, 1.909, 4, 3.1, -0.5, 1.975, 1.825, 2.75, 1.975, 1.825, 0.909, -1, -1, 0.9750000000000001, -1, 0.9750000000000001, -1)
$col = 11
foreach ($v in $row122vals) {
    $ws.Cells.Item(122, $col).Value = $v
    $col = $col + 1
}

# Row 123: JK Nomme United 0-0 JK Nomme Kalju
$ws.Cells.Item(123, 1).Value = 121
$ws.Cells.Item(123, 2).Value = 7723750
$ws.Cells.Item(123, 3).Value = "Estonia Meistriliiga"
$ws.Cells.Item(123, 4).Value = "Estonia Meistriliiga"
$ws.Cells.Item(123, 5).Value = 45368.39583333334
$ws.Cells.Item(123, 6).Value = "JK Nomme United"
$ws.Cells.Item(123, 7).Value = "JK Nomme Kalju"
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = "D"
$row123vals = @(5, 4, 1.5, 7, 4.5, 1.363, 1.5, 1.85, 1.95, 3, 1.775, 2.025, -1, 3.5, -1, 0.8500000000000001, -1, -1, 1.025)
$col = 11
foreach ($v in $row123vals) {
    $ws.Cells.Item(123, $col).Value = $v
    $col = $col + 1
}
